{"js": "// Add bright-green highlighting to the checklist's role/section labels\n// (Comments, Status, Photos, Diagnoses, Prescriptions - including their\n// \"(View only)\" variants) and to three stray notes at the bottom of the\n// doc (Location and alerts for animal?, 100 simultaneous users\n// (concurrency), All browsers).\n//\n// For the bulleted checklist items the highlight must land on BOTH the\n// paragraph mark (so the empty end-of-line pilcrow is highlighted too)\n// and the run itself - setting `paragraph.font.highlightColor` does\n// both in one shot. For the three plain paragraphs only the run text\n// itself is highlighted (the paragraph mark is left untouched), so\n// those are handled through a text search instead.\n\n// Paragraphs where the WHOLE paragraph text matches one of these labels\n// get the paragraph-level highlight (run + paragraph mark).\nconst paragraphTargets = new Set([\n  \"Comments\",\n  \"Status\",\n  \"Photos\",\n  \"Diagnoses\",\n  \"Prescriptions\",\n  \"Status (View only)\",\n  \"Photos (View only)\",\n  \"Diagnoses (View only)\",\n  \"Prescriptions (View only)\",\n]);\n\n// Plain paragraphs where only the run text (not the paragraph mark)\n// gets the highlight.\nconst runOnlyTargets = [\n  \"Location and alerts for animal?\",\n  \"100 simultaneous users (concurrency)\",\n  \"All browsers\",\n];\n\nconst HIGHLIGHT = \"BrightGreen\"; // -> OOXML <w:highlight w:val=\"green\"/>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (paragraphTargets.has(text)) {\n    paragraph.font.highlightColor = HIGHLIGHT;\n  }\n}\nawait context.sync();\n\nfor (const target of runOnlyTargets) {\n  const results = body.search(target, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n  for (const range of results.items) {\n    if (range.text === target) {\n      range.font.highlightColor = HIGHLIGHT;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Add bright-green highlighting to the checklist's role/section labels\n# (Comments, Status, Photos, Diagnoses, Prescriptions - including their\n# \"(View only)\" variants) and to three stray notes at the bottom of the\n# doc (Location and alerts for animal?, 100 simultaneous users\n# (concurrency), All browsers).\n#\n# For the bulleted checklist items the highlight must land on BOTH the\n# paragraph mark (so the empty end-of-line pilcrow is highlighted too)\n# and the run itself - setting the Range's Font.HighlightColorIndex\n# (rather than the Range's own HighlightColorIndex) achieves both in\n# one shot because Paragraph.Range includes the trailing paragraph\n# mark. For the three plain paragraphs only the run text itself is\n# highlighted (the paragraph mark is left untouched), so those are\n# handled through Find on just the text instead.\n\n$d = $word.ActiveDocument\n\n$paragraphTargets = @(\n    \"Comments\",\n    \"Status\",\n    \"Photos\",\n    \"Diagnoses\",\n    \"Prescriptions\",\n    \"Status (View only)\",\n    \"Photos (View only)\",\n    \"Diagnoses (View only)\",\n    \"Prescriptions (View only)\"\n)\n\n$runOnlyTargets = @(\n    \"Location and alerts for animal?\",\n    \"100 simultaneous users (concurrency)\",\n    \"All browsers\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paragraphTargets -contains $text) {\n        $p.Range.Font.HighlightColorIndex = \"wdBrightGreen\"\n    }\n}\n\nforeach ($target in $runOnlyTargets) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $target\n    $r.Find.MatchCase = $true\n    $r.Find.MatchWholeWord = $false\n    while ($r.Find.Execute()) {\n        if ($r.Text -eq $target) {\n            $r.HighlightColorIndex = \"wdBrightGreen\"\n        }\n        $r.Collapse(0)\n    }\n}\n"}
